# Use Case.docx - "Correções no modelo de domínio e to-do nos use cases"
#
# Rework the to-do block that follows the
# "(MOVIMENTO DE) SALDO USE CASES ..." heading near the end of the body.

$d = $word.ActiveDocument

# --- helpers ------------------------------------------------------------

# Append a brand-new paragraph right after paragraph number $n (1-based,
# per $d.Paragraphs.Item). $texts is an array of one or more text chunks;
# when there is more than one chunk each chunk ends up in its own <w:r>
# (built by creating temporary paragraph breaks and then splicing them
# back out, which keeps the runs distinct instead of merging them into a
# single <w:t>). The new paragraph inherits whatever direct paragraph
# formatting its predecessor ($n) currently has - callers that need a
# plain/default paragraph should only call this while the predecessor
# chain is still unformatted, and apply any w:pPr tweaks afterwards via
# Set-TightSpacing.
# Returns the paragraph number of the newly created paragraph.
function Add-ParaAfter($n, [string[]]$texts) {
    $prev = $d.Paragraphs.Item($n)
    $prev.Range.InsertParagraphAfter()
    $newN = $n + 1
    $newPara = $d.Paragraphs.Item($newN)
    $newPara.Range.Text = $texts[0]

    $curN = $newN
    for ($i = 1; $i -lt $texts.Length; $i++) {
        $curPara = $d.Paragraphs.Item($curN)
        $curPara.Range.InsertParagraphAfter()
        $chunkPara = $d.Paragraphs.Item($curN + 1)
        $chunkPara.Range.Text = $texts[$i]
        $curPara = $d.Paragraphs.Item($curN)
        $mark = $d.Range($curPara.Range.End - 1, $curPara.Range.End)
        $mark.Delete()
    }

    return $newN
}

# Give paragraph number $n the same "single spacing, auto" override the
# rest of the to-do list already uses (<w:spacing w:line="240"
# w:lineRule="auto"/>).
function Set-TightSpacing($n) {
    $para = $d.Paragraphs.Item($n)
    $para.Format.LineSpacingRule = 0   # wdLineSpaceSingle
    $para.Format.LineSpacing = 12      # -> w:line="240"
}

# --- 1. Locate the anchor paragraph ------------------------------------
$count = $d.Paragraphs.Count
$anchorN = -1
for ($i = $count; $i -ge 1; $i--) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*SALDO USE CASES*") {
        $anchorN = $i
        break
    }
}

$anchorPara = $d.Paragraphs.Item($anchorN)

# --- 2. Trim the heading text ------------------------------------------
$replaced = $anchorPara.Range.Find.Execute("SALDO USE CASES PRÓPRIOS E VER NOS OUTROS ONDE NECESSÁRIO", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "SALDO USE CASES PRÓPRIOS", 2)

# --- 3. Insert the new to-do paragraphs (still plain/default format) --
$n = $anchorN
$n = Add-ParaAfter $n @("Efetuar login – pós condição mudar?")
$n = Add-ParaAfter $n @("Criar Leilão – fazer a cena de retornar ao passo")
$n = Add-ParaAfter $n @("Adicionar ou Remover Leilões dos Favoritos", " – talvez dividir em dois???")
$n = Add-ParaAfter $n @("Efetuar Transação – tratar do saldo")
$n = Add-ParaAfter $n @("Pesquisar Leilões - ", "fazer a cena de retornar ao passo")
$n = Add-ParaAfter $n @("Consultar transações feitas - ", "fazer a cena de retornar ao passo")
$n = Add-ParaAfter $n @("Participar no leilão - ", "fazer a cena de retornar ao passo", " e tratar do saldo")
$n = Add-ParaAfter $n @("Nos use cases de tratar das contas, ver como fazer com o saldo")

# --- 4. Now that the whole chain exists, add the manual line break at
#        the end of the heading and apply the tight spacing only to the
#        heading + the first four to-do items -------------------------
$anchorPara = $d.Paragraphs.Item($anchorN)
$brPoint = $d.Range($anchorPara.Range.End - 1, $anchorPara.Range.End - 1)
$brPoint.InsertBreak(6)   # wdLineBreak -> <w:br/>

Set-TightSpacing $anchorN
Set-TightSpacing ($anchorN + 1)
Set-TightSpacing ($anchorN + 2)
Set-TightSpacing ($anchorN + 3)
Set-TightSpacing ($anchorN + 4)

# --- 5. Trim the trailing empty paragraphs down to a single one, then
#        turn the last of the six original trailing paragraphs into the
#        closing remark --------------------------------------------------
# $n is "Nos use cases..."; originally it was followed by 6 empty
# paragraphs. The very last paragraph of the body can never be deleted
# (Word always keeps a final paragraph mark), so keep paragraph $n+1
# empty, delete the next four (positions 2-5 of the six), and write the
# closing remark straight into the sixth (now final) paragraph instead
# of inserting a brand-new one after it.
for ($i = 0; $i -lt 4; $i++) {
    $victim = $d.Paragraphs.Item($n + 2)
    $victim.Range.Delete()
}

$closing = $d.Paragraphs.Item($n + 2)
$closing.Range.Text = "NÃO MEXI NO DIAGRAMA DOS USE CASES POIS SE CALHAR VAIS SEPARAR AQUELE E O DO SALDO SE CALHAR TMB SERIA BEM SEPARADO (podiamos pensar no extends nos do saldo quando tens de mexer nele, n sei como funciona muito bem)"

Write-Output "done"
